$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D10").Value = 0.1383952194513708
$ws.Range("D11").Value = 0.1383952194513708
$ws.Range("D12").Value = 0.000000000000000003090678116214709
$ws.Range("D13").Value = 0.000000000000000003090678116214709
